$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.202970297029703
$ws.Range("C2").Value = 0.5445544554455446
$ws.Range("J2").Value = 0.009900990099009901
$ws.Range("O2").Value = 0.004950495049504951
$ws.Range("P2").Value = 0.1683168316831683
$ws.Range("S2").Value = 0.06930693069306931
$ws.Range("B3").Value = 0.01739130434782609
$ws.Range("C3").Value = 0.04347826086956522
$ws.Range("J3").Value = 0.02608695652173913
$ws.Range("P3").Value = 0.7130434782608696
$ws.Range("S3").Value = 0.2
$ws.Range("J4").Value = 0.07692307692307693
$ws.Range("O4").Value = 0.03846153846153846
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.1923076923076923
$ws.Range("B6").Value = 0.06896551724137931
$ws.Range("D6").Value = 0.005747126436781609
$ws.Range("F6").Value = 0.06321839080459771
$ws.Range("J6").Value = 0.2413793103448276
$ws.Range("O6").Value = 0.01149425287356322
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.103448275862069
$ws.Range("S6").Value = 0.3390804597701149
$ws.Range("B7").Value = 0.1052631578947368
$ws.Range("D7").Value = 0.04093567251461988
$ws.Range("F7").Value = 0.02923976608187134
$ws.Range("J7").Value = 0.1228070175438596
$ws.Range("O7").Value = 0.005847953216374269
$ws.Range("Q7").Value = 0.1871345029239766
$ws.Range("R7").Value = 0.08187134502923976
$ws.Range("S7").Value = 0.4269005847953216
$ws.Range("B8").Value = 0.08455882352941177
$ws.Range("D8").Value = 0.01838235294117647
$ws.Range("F8").Value = 0.04779411764705882
$ws.Range("J8").Value = 0.1838235294117647
$ws.Range("O8").Value = 0.01838235294117647
$ws.Range("Q8").Value = 0.2169117647058824
$ws.Range("R8").Value = 0.09558823529411764
$ws.Range("S8").Value = 0.3345588235294117
$ws.Range("B9").Value = 0.08527131782945736
$ws.Range("F9").Value = 0.04651162790697674
$ws.Range("J9").Value = 0.124031007751938
$ws.Range("O9").Value = 0.03875968992248062
$ws.Range("Q9").Value = 0.1705426356589147
$ws.Range("R9").Value = 0.08527131782945736
$ws.Range("S9").Value = 0.4496124031007752
$ws.Range("B10").Value = 0.09547738693467336
$ws.Range("D10").Value = 0.01608040201005025
$ws.Range("F10").Value = 0.0814070351758794
$ws.Range("J10").Value = 0.1386934673366834
$ws.Range("O10").Value = 0.01206030150753769
$ws.Range("Q10").Value = 0.1798994974874372
$ws.Range("R10").Value = 0.09748743718592964
$ws.Range("S10").Value = 0.3788944723618091
$ws.Range("G11").Value = 0.1529850746268657
$ws.Range("J11").Value = 0.08582089552238806
$ws.Range("K11").Value = 0.2014925373134328
$ws.Range("L11").Value = 0.5335820895522388
$ws.Range("S11").Value = 0.02611940298507463
$ws.Range("G12").Value = 0.7162162162162162
$ws.Range("J12").Value = 0.1891891891891892
$ws.Range("K12").Value = 0.01351351351351351
$ws.Range("L12").Value = 0.03378378378378379
$ws.Range("S12").Value = 0.0472972972972973
$ws.Range("F13").Value = 0.02564102564102564
$ws.Range("G13").Value = 0.717948717948718
$ws.Range("J13").Value = 0.2051282051282051
$ws.Range("S13").Value = 0.05128205128205128
$ws.Range("F15").Value = 0.03289473684210526
$ws.Range("H15").Value = 0.1513157894736842
$ws.Range("I15").Value = 0.07236842105263158
$ws.Range("J15").Value = 0.3486842105263158
$ws.Range("K15").Value = 0.1052631578947368
$ws.Range("M15").Value = 0.0131578947368421
$ws.Range("O15").Value = 0.05921052631578947
$ws.Range("S15").Value = 0.2171052631578947
$ws.Range("F16").Value = 0.007751937984496124
$ws.Range("H16").Value = 0.1627906976744186
$ws.Range("I16").Value = 0.07751937984496124
$ws.Range("J16").Value = 0.3798449612403101
$ws.Range("K16").Value = 0.1472868217054264
$ws.Range("M16").Value = 0.04651162790697674
$ws.Range("O16").Value = 0.06201550387596899
$ws.Range("S16").Value = 0.1162790697674419
$ws.Range("F17").Value = 0.006309148264984227
$ws.Range("H17").Value = 0.1324921135646688
$ws.Range("I17").Value = 0.1041009463722398
$ws.Range("J17").Value = 0.4542586750788644
$ws.Range("K17").Value = 0.1072555205047319
$ws.Range("M17").Value = 0.01261829652996845
$ws.Range("O17").Value = 0.0473186119873817
$ws.Range("S17").Value = 0.1356466876971609
$ws.Range("F18").Value = 0.01212121212121212
$ws.Range("H18").Value = 0.1757575757575758
$ws.Range("I18").Value = 0.04242424242424243
$ws.Range("J18").Value = 0.4666666666666667
$ws.Range("K18").Value = 0.08484848484848485
$ws.Range("M18").Value = 0.04848484848484848
$ws.Range("N18").Value = 0.006060606060606061
$ws.Range("O18").Value = 0.07878787878787878
$ws.Range("S18").Value = 0.08484848484848485
$ws.Range("F19").Value = 0.01882845188284519
$ws.Range("H19").Value = 0.1642259414225941
$ws.Range("I19").Value = 0.07531380753138076
$ws.Range("J19").Value = 0.3692468619246862
$ws.Range("K19").Value = 0.1297071129707113
$ws.Range("M19").Value = 0.02092050209205021
$ws.Range("O19").Value = 0.06799163179916318
$ws.Range("S19").Value = 0.153765690376569
